# Commit: "Divided CSS file into multiple files"
# Adds a new "Uke 23" block (rows 31-32) to the plan, listing three follow-up
# tasks (improve assignment per feedback, image resolution, splitting/minifying
# the CSS file) with their Ja/Nei status markers, mirroring the layout of the
# existing weekly blocks further up the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New content: row 31 ---------------------------------------------------
$ws.Range("A31").Value = "Uke 23"
$ws.Range("B31").Value = "Forbedre oppgaven etter tilbakemeldningen"
$ws.Range("C31").Value = "Nei"

# --- New content: row 32 ----------------------------------------------------
# (set before D31 so the new shared strings are appended in the same order
#  as in the target workbook)
$ws.Range("B32").Value = " Oppløsning på bilder"

$ws.Range("D31").Value = "Forminke CSS filen"
$ws.Range("E31").Value = "Nei"
$ws.Range("C32").Value = "Ja"

# --- Formatting: reuse the same styles used by the other weekly blocks -----
# A-column week header (bold)
$ws.Range("A27").Copy()
$ws.Range("A31").PasteSpecial(-4122)

# "Nei" markers (red-ish fill)
$ws.Range("C28").Copy()
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("E31").PasteSpecial(-4122)

# "Ja" marker (green-ish fill)
$ws.Range("C27").Copy()
$ws.Range("C32").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- View state: selection moved down near the newly added rows ------------
$ws.Range("F33").Select()
